$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) updates ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 51
$ws1.Range("F3").Value = 21472
$ws1.Range("G3").Value = 70
$ws1.Range("F4").Value = 819
$ws1.Range("F6").Value = 1135
$ws1.Range("F8").Value = 7993
$ws1.Range("F9").Value = 558
$ws1.Range("F10").Value = 45
$ws1.Range("F13").Value = 69
$ws1.Range("F15").Value = 180
$ws1.Range("F16").Value = 35
$ws1.Range("F20").Value = 548
$ws1.Range("F22").Value = 711
$ws1.Range("F24").Value = 87
$ws1.Range("F27").Value = 1201
$ws1.Range("F28").Value = 58
$ws1.Range("F29").Value = 43
$ws1.Range("F30").Value = 230
$ws1.Range("F31").Value = 607
$ws1.Range("F33").Value = 145
$ws1.Range("F34").Value = 5102
$ws1.Range("F37").Value = 53
$ws1.Range("F38").Value = 59
$ws1.Range("F39").Value = 13206
$ws1.Range("F40").Value = 1370
$ws1.Range("F41").Value = 144
$ws1.Range("F42").Value = 58
$ws1.Range("F44").Value = 318
$ws1.Range("F46").Value = 4071
$ws1.Range("F47").Value = 3
$ws1.Range("F48").Value = 333

# --- Sheet "全部类型" (sheet4) updates ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 51
$ws4.Range("F3").Value = 21472
$ws4.Range("G3").Value = 70
$ws4.Range("F4").Value = 1135
$ws4.Range("F6").Value = 7993
$ws4.Range("F7").Value = 558
$ws4.Range("F8").Value = 45
$ws4.Range("F11").Value = 69
$ws4.Range("F13").Value = 180
$ws4.Range("F14").Value = 35
$ws4.Range("F17").Value = 548
$ws4.Range("F19").Value = 711
$ws4.Range("F21").Value = 87
$ws4.Range("F24").Value = 1201
$ws4.Range("F25").Value = 58
$ws4.Range("F26").Value = 43
$ws4.Range("F27").Value = 230
$ws4.Range("F29").Value = 607
$ws4.Range("F32").Value = 145
$ws4.Range("F34").Value = 5102
$ws4.Range("F37").Value = 53
$ws4.Range("F38").Value = 59
$ws4.Range("F39").Value = 13207
$ws4.Range("F40").Value = 1370
$ws4.Range("F41").Value = 144
$ws4.Range("F42").Value = 58
$ws4.Range("F44").Value = 318
$ws4.Range("F46").Value = 4071
$ws4.Range("F47").Value = 3
$ws4.Range("F48").Value = 333

